# Applies the "Current Submission" sheet update:
#  - Row 3 header labels keep their text (Tag/Date/Submission/Description of
#    Changes) - untouched so the shared-string table keeps its existing slots
#  - A1 header text becomes two-run rich text (plain prefix + bold team/org
#    name suffix)
#  - Row 4 (previously the "TurnIn1" entry) is replaced with the new
#    "ActDue01 / Team Project Ideas" submission
#  - Four more submission rows (5-8) are added for ActDue02..ActDue05
#  - Row heights bump slightly (62 -> 62.1) to match the re-saved workbook
#  - Selection moves from A5 to B9

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Rows 4-8: submission entries. Values are written in the same order the
# original author entered them (tags for every row, then content filled in
# out of row order) so the shared-string table lines up with the source.
# ---------------------------------------------------------------------
$ws.Cells.Item(4, 1).Value = "ActDue01"
$ws.Cells.Item(4, 3).Value = "Team Project Ideas"
$ws.Cells.Item(4, 4).Value = "Team formation, proposed projects, summary of work to be done."

$ws.Cells.Item(5, 1).Value = "ActDue02"
$ws.Cells.Item(6, 1).Value = "ActDue03"
$ws.Cells.Item(7, 1).Value = "ActDue04"
$ws.Cells.Item(8, 1).Value = "ActDue05"

$ws.Cells.Item(7, 3).Value = "Software Selection"

$ws.Cells.Item(5, 3).Value = "Software Requirements"
$ws.Cells.Item(5, 4).Value = "Software Requirements Specification - Use cases division"

$ws.Cells.Item(6, 3).Value = "Finshed SRS and posted"
$ws.Cells.Item(6, 4).Value = "Finished SRS and posted to the BlackBoard"

$ws.Cells.Item(7, 4).Value = "Selected Java and started creating the project classes and screens"

$ws.Cells.Item(8, 3).Value = "First Sprint Started"
$ws.Cells.Item(8, 4).Value = "Programminng Frames : Eric(SetupSummary), Ruth (AnswerEntry), Brian (QuestionSelect), Laercio (GameStatus) and Alesander (WinnerSelect)"

# Dates (column B) - numeric, not part of the shared-string table
$ws.Cells.Item(4, 2).Value = 40789
$ws.Cells.Item(5, 2).Value = 40810
$ws.Cells.Item(6, 2).Value = 40817
$ws.Cells.Item(7, 2).Value = 40824
$ws.Cells.Item(8, 2).Value = 40831

# ---------------------------------------------------------------------
# Formatting: wrap text on the Submission/Description columns (row 6's
# Description column is the one exception that stays unwrapped), font for
# the tag cells (and the empty row9 tag cell), and row heights
# ---------------------------------------------------------------------
foreach ($r in 4..9) {
  $aCell = $ws.Cells.Item($r, 1)
  $aCell.Font.Name = "Verdana"
  $aCell.Font.Size = 10
  $aCell.Font.Bold = $false
}

foreach ($r in 4..8) {
  $ws.Cells.Item($r, 3).WrapText = $true
  $ws.Cells.Item($r, 4).WrapText = $true

  $ws.Rows.Item($r).RowHeight = 62.1
}
$ws.Cells.Item(6, 4).WrapText = $false

foreach ($addr in @("C4", "D4")) {
  $c = $ws.Range($addr)
  $c.Font.Name = "Verdana"
  $c.Font.Size = 10
  $c.Font.Bold = $false
}

# ---------------------------------------------------------------------
# Normalize row height for the remaining (still-empty) rows 9-25 to match
# the re-saved file's slightly larger default height
# ---------------------------------------------------------------------
for ($n = 9; $n -le 25; $n++) {
  $ws.Rows.Item($n).RowHeight = 62.1
}

# ---------------------------------------------------------------------
# A1: rich-text header - plain prefix + bold team/org name. Done after the
# other string writes since rich-text runs are always serialized last.
# ---------------------------------------------------------------------
$prefix = "Current Submission for Team ___"
$suffix = "Organization of Programming Students, OPS"
$ws.Range("A1").Value = $prefix + $suffix
$boldChars = $ws.Range("A1").Characters($prefix.Length + 1, $suffix.Length)
$boldChars.Font.Bold = $true
$boldChars.Font.Name = "Verdana"
$boldChars.Font.Size = 10

# ---------------------------------------------------------------------
# Selection moves to B9 (first empty submission row) in the saved file
# ---------------------------------------------------------------------
$ws.Range("B9").Select()
